# Handback report generation: refresh the handoff/handback timestamps
# recorded against the 7a82f988-... file after a new round-trip, and roll
# the newest of those timestamps up into the Overview sheet's
# "Latest HO Xliff Generate Date" column for that row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K)
# for the 7a82f988-... row (row 3).
$wsZhCn.Range("H3").Value = "2016-08-22 02:58:08"
$wsZhCn.Range("K3").Value = "2016-08-22 02:58:25"

# de-de: same two columns, same row.
$wsDeDe.Range("H3").Value = "2016-08-22 02:58:13"
$wsDeDe.Range("K3").Value = "2016-08-22 02:58:32"

# Overview: "Latest HO Xliff Generate Date" for that file is the max of the
# per-language handoff datetimes just written above.
$wsOverview.Range("G3").Value = "2016-08-22 02:58:13"
